# "Requisitos Funcionais" workbook update.
#
# The requirement row "RF07 - Gerar Relatório de Alimentação" (a duplicate /
# superseded requirement, whose description duplicated the later "RF08 -
# Consultar Frequência de Alimentação" theme) was removed from the
# requirements table on sheet "Planilha1". Deleting the entire worksheet row
# shifts every subsequent row up by one; the "Código" column then needs its
# sequential RFxx labels re-numbered so the codes stay contiguous
# (old RF08..RF17 -> new RF07..RF16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RF07 - Gerar Relatório de Alimentação" row (row 8: Código,
# Nome, Descrição, Prioridade, Categoria). Everything below shifts up.
$ws.Rows(8).Delete()

# Re-sequence the "Código" column for the rows that shifted up so the
# RFxx numbering stays contiguous (no gap / no duplicate).
$ws.Range("A8").Value = "RF07"
$ws.Range("A9").Value = "RF08"
$ws.Range("A10").Value = "RF09"
$ws.Range("A11").Value = "RF10"
$ws.Range("A12").Value = "RF11"
$ws.Range("A13").Value = "RF12"
$ws.Range("A14").Value = "RF13"
$ws.Range("A15").Value = "RF14"
$ws.Range("A16").Value = "RF15"
$ws.Range("A17").Value = "RF16"

# Match the author's final cursor position after the edit.
$ws.Range("B19").Select()
